# Scheduled-runner refresh of market-board figures on the "Golem_Profits"
# workbook. Re-prices currentAveragePrice* (H/I/J) and the downstream
# LevePrice*/LeveProfit* (K/L/M/N) columns for the leves whose market
# snapshot changed since the last run; blank cells are cleared where the
# refreshed row no longer carries a NQ/HQ split.
$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")
# Row 46: Always Have an Exit Plan / Poisoning Potion
$ws.Range("H46").Value = 1166.6666
$ws.Range("J46").Value = 750
$ws.Range("L46").Value = 2250
$ws.Range("N46").Value = -2488
# Row 60: Make Up Your Mind or Else / Potent Poisoning Potion
$ws.Range("H60").Value = 1166.6666
$ws.Range("J60").Value = 750
$ws.Range("L60").Value = 2250
$ws.Range("N60").Value = -3218
# Row 112: Making Ends Meet / Superior Spiritbond Potion
$ws.Range("H112").Value = 1949.9
$ws.Range("J112").Value = 1999.8889
$ws.Range("L112").Value = 5999.6667
$ws.Range("N112").Value = -8215.6667
# Row 140: Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").Value = ""

# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")
# Row 6: Don't Hit Me One More Time / Bronze Hoplon
$ws.Range("H6").Value = 502
$ws.Range("I6").Value = 502
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 502
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -329
$ws.Range("N6").Value = ""
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3395.2
$ws.Range("I61").Value = 3395.2
$ws.Range("K61").Value = 3395.2
$ws.Range("M61").Value = -3183.2
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2122.5
$ws.Range("I74").Value = 2122.5
$ws.Range("K74").Value = 2122.5
$ws.Range("M74").Value = -1248.5
# Row 76: Sometimes the South Wins / Titanium Mail of Fending
$ws.Range("H76").Value = 28198.334
$ws.Range("J76").Value = 28198.334
$ws.Range("L76").Value = 28198.334
$ws.Range("N76").Value = -28874.334
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2122.5
$ws.Range("I77").Value = 2122.5
$ws.Range("K77").Value = 10612.5
$ws.Range("M77").Value = -6244.5
# Row 79: The Thriller of Autumn (L) / Titanium Mail of Fending
$ws.Range("H79").Value = 28198.334
$ws.Range("J79").Value = 28198.334
$ws.Range("L79").Value = 28198.334
$ws.Range("N79").Value = -30538.334
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 1550.6666
$ws.Range("I122").Value = 1550.6666
$ws.Range("K122").Value = 4651.9998
$ws.Range("M122").Value = -2201.9998
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1770.909
$ws.Range("I132").Value = 1748
$ws.Range("K132").Value = 5244
$ws.Range("M132").Value = -2714
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3395.2
$ws.Range("I136").Value = 3395.2
$ws.Range("K136").Value = 10185.6
$ws.Range("M136").Value = -7635.599999999999
# Row 137: Odd Instruments / Cobalt Tungsten Alembic
$ws.Range("H137").Value = 5000
$ws.Range("I137").Value = 5000
$ws.Range("K137").Value = 5000
$ws.Range("M137").Value = 100
# Row 140: A Hand for a Deckhand / Ra'Kaznar Gloves of Scouting
$ws.Range("H140").Value = 50429
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 50429
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 50429
$ws.Range("M140").Value = ""
$ws.Range("N140").Value = -60789
# Row 141: Essays on Equipment / Ra'Kaznar Greaves of Maiming
$ws.Range("H141").Value = 32500
$ws.Range("I141").Value = 5000
$ws.Range("J141").Value = 60000
$ws.Range("K141").Value = 5000
$ws.Range("L141").Value = 60000
$ws.Range("M141").Value = 180
$ws.Range("N141").Value = -70360

# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 1252.8
$ws.Range("J20").Value = 1386
$ws.Range("L20").Value = 1386
$ws.Range("N20").Value = -1880

# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 194.14285
$ws.Range("I7").Value = 191.25
$ws.Range("J7").Value = 198
$ws.Range("K7").Value = 191.25
$ws.Range("L7").Value = 198
$ws.Range("M7").Value = -78.25
$ws.Range("N7").Value = -424

# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")
# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 1916.5
$ws.Range("I68").Value = 1854.4
$ws.Range("J68").Value = 1944.7273
$ws.Range("K68").Value = 5563.200000000001
$ws.Range("L68").Value = 5834.1819
$ws.Range("M68").Value = -4752.200000000001
$ws.Range("N68").Value = -7456.1819
# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 1916.5
$ws.Range("I71").Value = 1854.4
$ws.Range("J71").Value = 1944.7273
$ws.Range("K71").Value = 16689.6
$ws.Range("L71").Value = 17502.5457
$ws.Range("M71").Value = -12633.6
$ws.Range("N71").Value = -25614.5457
# Row 137: Creative Chocolate / Gateau au Chocolat
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = ""

# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 125002760
$ws.Range("I70").Value = 3275
$ws.Range("J70").Value = 250002260
$ws.Range("K70").Value = 3275
$ws.Range("L70").Value = 250002260
$ws.Range("M70").Value = -3005
$ws.Range("N70").Value = -250002800
# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 125002760
$ws.Range("I73").Value = 3275
$ws.Range("J73").Value = 250002260
$ws.Range("K73").Value = 3275
$ws.Range("L73").Value = 250002260
$ws.Range("M73").Value = -2339
$ws.Range("N73").Value = -250004132

# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1330
$ws.Range("N16").Value = ""
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""
# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""
# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 1457
$ws.Range("I46").Value = 833.3333
$ws.Range("J46").Value = 2392.5
$ws.Range("K46").Value = 833.3333
$ws.Range("L46").Value = 2392.5
$ws.Range("M46").Value = -645.3333
$ws.Range("N46").Value = -2768.5
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 4949.8
$ws.Range("I136").Value = 4249.25
$ws.Range("K136").Value = 12747.75
$ws.Range("M136").Value = -10197.75

# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")
# Row 70: An Account of My Boots / Holy Rainbow Shoes
$ws.Range("H70").Value = 60000
$ws.Range("I70").Value = 60000
$ws.Range("K70").Value = 60000
$ws.Range("M70").Value = -59685
# Row 73: Soot in My Hair and Scars on My Feet (L) / Holy Rainbow Shoes
$ws.Range("H73").Value = 60000
$ws.Range("I73").Value = 60000
$ws.Range("K73").Value = 60000
$ws.Range("M73").Value = -58908
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1987.5
$ws.Range("I122").Value = 1650
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4950
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2500
$ws.Range("N122").Value = -13900
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1918
$ws.Range("I136").Value = 1918
$ws.Range("K136").Value = 5754
$ws.Range("M136").Value = -3204
